$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jurisdiction")
$ws.Range("B6").Value = 100
Write-Output "B6 type: $($ws.Range('B6').Value2.GetType())"
